$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) and G (Hora) hold plain-text numeric /
# percentage strings in the source data (Excel would otherwise auto-convert
# values like "307.52" or "19" to Number on assignment). Force Text format on
# each touched cell immediately before writing its new value so the stored
# type/string representation is preserved exactly as in the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.52"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.98%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "19"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.66"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.47%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "19"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.031"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.71%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "19"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07878"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.57%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "19"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.177"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.62%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "19"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.060"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.76%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "19"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.062"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.47%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "19"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9258"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.05%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "19"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09936"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.87%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "19"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1878"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.06%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "19"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08707"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.07%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "19"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03616"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.39%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "19"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09945"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.12%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "19"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001482"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.22%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "19"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005655"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.34%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "19"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.462"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.52%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "19"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.469"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "16.10%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "19"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3452"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.57%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "19"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1337"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.27%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "19"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.931"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "8.39%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "19"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2202"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.49%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "19"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04622"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.98%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "19"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.005205"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "15.53%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "19"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001233"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.49%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "19"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001401"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "7.79%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "19"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002720"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "19"

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "19"

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "19"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "19"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "19"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "19"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "19"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "19"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "19"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "19"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "19"

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "19"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01822"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.35%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "19"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04749"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.92%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "19"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007902"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.39%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "19"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1409"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.90%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "19"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007606"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.49%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "19"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002182"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.72%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "19"

$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "19"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006317"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.59%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "19"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.14%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "19"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0005806"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.10%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "19"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.65"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "796.11%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "19"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.15%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "19"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00002102"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.14%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "19"
